$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

# Copy formatting (style) from the last existing data row (78) so the
# new row matches the table's look, then overwrite the values.
$ws.Range("A78:F78").Copy()
$ws.Range("A79:F79").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value() = "L_ERSTT_4"
$ws.Cells.Item($row, 2).Value() = "Q_EUROSTAT"
$ws.Cells.Item($row, 3).Value() = "Rate der (erheblichen) materiellen Deprivation"
$ws.Cells.Item($row, 4).Value() = "(Severe) Material deprivation rate"
$ws.Cells.Item($row, 5).Value() = "https://ec.europa.eu/eurostat/databrowser/view/ILC_SIP8__custom_5385803/default/table?lang=de"
$ws.Cells.Item($row, 6).Value() = "https://ec.europa.eu/eurostat/databrowser/view/ILC_SIP8__custom_5385803/default/table?lang=en"
